$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handoff file being reported on: d08eb7dc-77a2-4176-a7f0-6948c41606dc.md
# Appears as a new row (row 9) on all three sheets: Overview, zh-cn, de-de.
# ---------------------------------------------------------------------------

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = "d08eb7dc-77a2-4176-a7f0-6948c41606dc.md"
$wsOverview.Range("B9").Value = "e2e\d08eb7dc-77a2-4176-a7f0-6948c41606dc.md"
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-20 10:49:42"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d08eb7dc77a24176a7f06948c41606dc0000000/e2e/d08eb7dc-77a2-4176-a7f0-6948c41606dc.md", "", "", "e2e\d08eb7dc-77a2-4176-a7f0-6948c41606dc.md") | Out-Null
$wsOverview.Range("B9").Style = "HyperLink"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A9").Value = "d08eb7dc-77a2-4176-a7f0-6948c41606dc.md"
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "d08eb7dc-77a2-4176-a7f0-6948c41606dc.235bfabee2edd399bafeff1ea485ddaf03dafd97.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-20 10:49:38"
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "False"
$wsZhCn.Range("P9").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d08eb7dc77a24176a7f06948c41606dc0000001/e2e/d08eb7dc-77a2-4176-a7f0-6948c41606dc.md", "", "", "d08eb7dc-77a2-4176-a7f0-6948c41606dc.md") | Out-Null
$wsZhCn.Range("A9").Style = "HyperLink"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A9").Value = "d08eb7dc-77a2-4176-a7f0-6948c41606dc.md"
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "d08eb7dc-77a2-4176-a7f0-6948c41606dc.235bfabee2edd399bafeff1ea485ddaf03dafd97.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-20 10:49:42"
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "False"
$wsDeDe.Range("P9").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d08eb7dc77a24176a7f06948c41606dc0000002/e2e/d08eb7dc-77a2-4176-a7f0-6948c41606dc.md", "", "", "d08eb7dc-77a2-4176-a7f0-6948c41606dc.md") | Out-Null
$wsDeDe.Range("A9").Style = "HyperLink"
